$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the price/volume columns so that numeric-looking
# strings (e.g. "218.79", "4.482") are written as text, matching the
# original inline-string cell contents, not auto-converted to numbers.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "26.616.76"
$ws.Range("E2").Value = "  -7.42%  "
$ws.Range("D3").Value = "1.696.46"
$ws.Range("E3").Value = "  -6.06%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "218.79"
$ws.Range("D6").Value = "0.5086"
$ws.Range("E6").Value = "  -14.03%  "
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "0.2645"
$ws.Range("E8").Value = "  -4.80%  "
$ws.Range("D9").Value = "22.17"
$ws.Range("E9").Value = "  -5.12%  "
$ws.Range("D10").Value = "0.06244"
$ws.Range("E10").Value = "  -8.63%  "
$ws.Range("D11").Value = "0.07314"
$ws.Range("E11").Value = "  -2.47%  "
$ws.Range("D12").Value = "1.697.38"
$ws.Range("E12").Value = "  -6.00%  "
$ws.Range("D13").Value = "4.482"
$ws.Range("E13").Value = "  -6.06%  "
$ws.Range("D14").Value = "0.5834"
$ws.Range("E14").Value = "  -6.56%  "
$ws.Range("D15").Value = "1.927.28"
$ws.Range("E15").Value = "  -6.05%  "
$ws.Range("D16").Value = "0.000008322"
$ws.Range("E16").Value = "  -10.45%  "
$ws.Range("D17").Value = "65.38"
$ws.Range("E17").Value = "  -13.78%  "
$ws.Range("D18").Value = "26.671.94"
$ws.Range("E18").Value = "  -7.03%  "
$ws.Range("D19").Value = "5.047"
$ws.Range("E19").Value = "  -8.00%  "
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").Value = "10.86"
$ws.Range("E21").Value = "  -5.46%  "
$ws.Range("D22").Value = "186.52"
$ws.Range("D23").Value = "6.258"
$ws.Range("E23").Value = "  -8.55%  "
$ws.Range("D24").Value = "1.006"
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("D25").Value = "144.94"
$ws.Range("E25").Value = "  -6.08%  "
$ws.Range("D26").Value = "7.607"
$ws.Range("E26").Value = "  -3.59%  "
$ws.Range("E27").Value = "  -10.01%  "
$ws.Range("E28").Value = "  -4.70%  "
$ws.Range("D29").Value = "1.298"
$ws.Range("E29").Value = "  -8.86%  "
$ws.Range("D30").Value = "0.05705"
$ws.Range("E30").Value = "  -7.66%  "
$ws.Range("D31").Value = "1.335"
$ws.Range("E31").Value = "  -6.33%  "
$ws.Range("D32").Value = "3.512"
$ws.Range("E32").Value = "  -6.69%  "
$ws.Range("D33").Value = "3.490"
$ws.Range("E33").Value = "  -7.79%  "
$ws.Range("D34").Value = "1.658"
$ws.Range("E34").Value = "  -4.27%  "
$ws.Range("E35").Value = "  -4.52%  "
$ws.Range("D36").Value = "0.5993"
$ws.Range("E36").Value = "  -6.75%  "
$ws.Range("D37").Value = "2.375"
$ws.Range("E37").Value = "  -4.75%  "
$ws.Range("D38").Value = "2.679"
$ws.Range("E38").Value = "  -1.68%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "1.088.81"
$ws.Range("E39").Value = "  -4.71%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "0.01595"
$ws.Range("E40").Value = "  -6.88%  "
$ws.Range("D41").Value = "5.910"
$ws.Range("E41").Value = "  -10.00%  "
$ws.Range("D42").Value = "0.8642"
$ws.Range("E42").Value = "  -1.92%  "
$ws.Range("E43").Value = "  -0.46%  "
$ws.Range("D44").Value = "98.58"
$ws.Range("E44").Value = "  -1.76%  "
$ws.Range("D45").Value = "1.854.39"
$ws.Range("E45").Value = "  -5.48%  "
$ws.Range("D46").Value = "56.68"
$ws.Range("E46").Value = "  -6.36%  "
$ws.Range("D47").Value = "0.00000000106"
$ws.Range("E47").Value = "  -4.84%  "
$ws.Range("D48").Value = "1.004"
$ws.Range("E48").Value = "  -0.72%  "
$ws.Range("D49").Value = "8.132"
$ws.Range("E49").Value = "  -2.71%  "
$ws.Range("D50").Value = "0.05242"
$ws.Range("E50").Value = "  -4.12%  "
$ws.Range("D51").Value = "0.4317"

# Restore the default (Normal) style so the cells keep the same
# unstyled appearance as before (no explicit style index).
$dataRange.Style = "Normal"
